# Migrando de Linq to SQL para Entity Framework
# Adds a new "PointsPerDay" entity field to the "Pessoa" sheet, removing the
# old free-text "Pontos por dia" note (the field it formalizes) and
# reshuffling the trailing "Campos a se pensar:" notes block accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pessoa")

# --- New formal field row: PointsPerDay ------------------------------------
$ws.Range("A9").Value = "PointsPerDay"
$ws.Range("B9").Value = "Pontos Por dia"
$ws.Range("C9").Value = "Informa quantos pontos a pessoa pode comer por dia, em alimentos."
$ws.Range("D9").Value = "Int16"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = "Não"

# --- "Campos a se pensar:" notes block: drop the now-formalised items and
#     move the section header down so the remaining note lines up under it --
$ws.Range("B11").Value = ""
$ws.Range("B11").Font.Bold = $false

$ws.Range("B12").Value = ""

$ws.Range("B13").Value = ""

$ws.Range("B14").Value = "Campos a se pensar:"
$ws.Range("B14").Font.Bold = $true

$ws.Range("B15").Value = "Pontos extras por semana"

# --- Widen the Observações column so the new long description fits --------
$ws.Columns.Item(3).ColumnWidth = 62.5

# --- Match the author's final selection ------------------------------------
[void]$ws.Range("B14").Select()
